$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '240.55'
$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = '1'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.38'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '1'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.487'
$ws.Range('G4').NumberFormat = '@'
$ws.Range('G4').Value = '1'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05580'
$ws.Range('G5').NumberFormat = '@'
$ws.Range('G5').Value = '1'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.372'
$ws.Range('E6').Value = '5GateTokenGT'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '1'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.468'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('G7').NumberFormat = '@'
$ws.Range('G7').Value = '1'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.076'
$ws.Range('E8').Value = '7FTXTokenFTT'
$ws.Range('G8').NumberFormat = '@'
$ws.Range('G8').Value = '1'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8032'
$ws.Range('E9').Value = '8MXTokenMX'
$ws.Range('G9').NumberFormat = '@'
$ws.Range('G9').Value = '1'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1422'
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value = '1'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07498'
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value = '1'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03257'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value = '1'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02983'
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value = '1'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09248'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value = '1'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001658'
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value = '1'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.249'
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value = '1'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04748'
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value = '1'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0005750'
$ws.Range('G18').NumberFormat = '@'
$ws.Range('G18').Value = '1'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006267'
$ws.Range('G19').NumberFormat = '@'
$ws.Range('G19').Value = '1'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.001050'
$ws.Range('G20').NumberFormat = '@'
$ws.Range('G20').Value = '1'
$ws.Range('G21').NumberFormat = '@'
$ws.Range('G21').Value = '1'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001500'
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '1'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0004777'
$ws.Range('G23').NumberFormat = '@'
$ws.Range('G23').Value = '1'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.973'
$ws.Range('E24').Value = '23LEOLEOBestin24h'
$ws.Range('G24').NumberFormat = '@'
$ws.Range('G24').Value = '1'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.142'
$ws.Range('G25').NumberFormat = '@'
$ws.Range('G25').Value = '1'
$ws.Range('G26').NumberFormat = '@'
$ws.Range('G26').Value = '1'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1313'
$ws.Range('G27').NumberFormat = '@'
$ws.Range('G27').Value = '1'
$ws.Range('G28').NumberFormat = '@'
$ws.Range('G28').Value = '1'
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '1'
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '1'
$ws.Range('G31').NumberFormat = '@'
$ws.Range('G31').Value = '1'
$ws.Range('G32').NumberFormat = '@'
$ws.Range('G32').Value = '1'
$ws.Range('G33').NumberFormat = '@'
$ws.Range('G33').Value = '1'
$ws.Range('G34').NumberFormat = '@'
$ws.Range('G34').Value = '1'
$ws.Range('G35').NumberFormat = '@'
$ws.Range('G35').Value = '1'
$ws.Range('G36').NumberFormat = '@'
$ws.Range('G36').Value = '1'
$ws.Range('G37').NumberFormat = '@'
$ws.Range('G37').Value = '1'
$ws.Range('G38').NumberFormat = '@'
$ws.Range('G38').Value = '1'
$ws.Range('G39').NumberFormat = '@'
$ws.Range('G39').Value = '1'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04164'
$ws.Range('G40').NumberFormat = '@'
$ws.Range('G40').Value = '1'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006996'
$ws.Range('G41').NumberFormat = '@'
$ws.Range('G41').Value = '1'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1043'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').NumberFormat = '@'
$ws.Range('G42').Value = '1'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002970'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('G43').NumberFormat = '@'
$ws.Range('G43').Value = '1'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008806'
$ws.Range('G44').NumberFormat = '@'
$ws.Range('G44').Value = '1'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005466'
$ws.Range('G45').NumberFormat = '@'
$ws.Range('G45').Value = '1'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('G46').NumberFormat = '@'
$ws.Range('G46').Value = '1'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6800'
$ws.Range('G47').NumberFormat = '@'
$ws.Range('G47').Value = '1'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.03066'
$ws.Range('G48').NumberFormat = '@'
$ws.Range('G48').Value = '1'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('G49').NumberFormat = '@'
$ws.Range('G49').Value = '1'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.01010'
$ws.Range('G50').NumberFormat = '@'
$ws.Range('G50').Value = '1'
$ws.Range('G51').NumberFormat = '@'
$ws.Range('G51').Value = '1'
